$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.685.01"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "1.803.58"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "'231.98"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").Value = "'0.5922"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "'0.2772"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "'0.06817"
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").Value = "'0.07509"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").Value = "1.803.56"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "'4.765"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "'0.6215"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "2.048.83"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "'0.000009213"
$ws.Range("E16").Value = "  -7.28%  "
$ws.Range("D17").Value = "'75.52"
$ws.Range("D18").Value = "28.644.67"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "'5.485"
$ws.Range("E19").Value = "  -6.21%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'210.64"
$ws.Range("E21").Value = "  -6.85%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "'6.827"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D25").Value = "'153.82"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "'7.861"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").Value = "'0.1268"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "'16.43"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "'1.424"
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("D30").Value = "'0.06196"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "'1.427"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").Value = "'3.781"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "'3.742"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").Value = "'1.733"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D36").Value = "'0.6426"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "'2.498"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").Value = "'2.714"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "'6.533"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "'0.01690"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "1.147.16"
$ws.Range("E41").Value = "  -5.74%  "
$ws.Range("D42").Value = "'0.8836"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "'1.005"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").Value = "'99.93"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "1.952.44"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("D48").Value = "'1.589"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'8.350"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").Value = "'0.05475"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'0.4478"
$ws.Range("E51").Value = "  -1.44%  "
